$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.465.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.72%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.465.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.93%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.87%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.465.68"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.87%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.485"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.53%  "

$ws.Range("E10").Value = "  -1.88%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.15"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("E12").Value = "  -2.71%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.061.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.69%  "

$ws.Range("E15").Value = "  +1.34%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000176"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.81%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.503.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.651.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.61%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "380.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.569"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.607.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.95%  "

$ws.Range("E25").Value = "  +0.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "72.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.29%  "

$ws.Range("E27").Value = "  -4.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.16%  "

$ws.Range("E31").Value = "  -3.74%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.472.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.94%  "

$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("E35").Value = "  -2.78%  "

$ws.Range("E36").Value = "  -2.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.34%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.54"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.02%  "

$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "166.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.59%  "

$ws.Range("E40").Value = "  -1.11%  "

$ws.Range("E41").Value = "  -4.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.807"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.38%  "

$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("E45").Value = "  -4.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.09%  "

$ws.Range("E48").Value = "  -3.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.415.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.70%  "

$ws.Range("E51").Value = "  -2.17%  "
